$d = $word.ActiveDocument

# Helper: insert a brand-new paragraph right after the paragraph currently
# at index $paraIndex, fill it with the given ordered list of runs
# (@{Text=...; Bold=$true/$false} hashtables), apply bold to the
# sub-ranges that need it, and return the index of the paragraph just
# created (so the caller can chain further insertions after it).
function Add-RunsParagraphAtIndex($paraIndex, $runs, $d) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.InsertParagraphAfter()

    $newIndex = $paraIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newRange = $newPara.Range
    $paraStart = $newRange.Start

    $fullText = ""
    foreach ($run in $runs) {
        $fullText = $fullText + $run.Text
    }
    $newRange.InsertAfter($fullText)

    $cursor = $paraStart
    foreach ($run in $runs) {
        $len = $run.Text.Length
        if ($run.Bold -and $len -gt 0) {
            $sub = $d.Range($cursor, $cursor + $len)
            $sub.Font.Bold = 1
            $sub.Font.BoldBi = 1
        }
        $cursor = $cursor + $len
    }

    return $newIndex
}

$dash = [string][char]0x2013
$ldq = [string][char]0x201C
$rdq = [string][char]0x201D

# Locate the paragraph that currently ends with "...in the intro scene."
# (the last populated paragraph before the two trailing blank paragraphs).
$idx = 14

# --- "Instructions Scene" heading (bold, 18pt / sz 36) -------------------
$p = $d.Paragraphs.Item($idx)
$rng = $p.Range
$rng.InsertParagraphAfter()
$idx = $idx + 1
$headingPara = $d.Paragraphs.Item($idx)
$headingRange = $headingPara.Range
$headingRange.InsertAfter("Instructions Scene")
$headingFull = $headingPara.Range
$headingFull.Font.Bold = 1
$headingFull.Font.BoldBi = 1
$headingFull.Font.Size = 18
$headingFull.Font.SizeBi = 18

# --- "Four visual elements:" ---------------------------------------------
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="Four visual element"; Bold=$false},
    @{Text="s:"; Bold=$false}
) $d

# --- "Instructions – SimpleGE multilabel explaining the game play instructions"
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="Instructions "; Bold=$true},
    @{Text=$dash; Bold=$true},
    @{Text=" "; Bold=$true},
    @{Text="SimpleGE"; Bold=$false},
    @{Text=" multilabel "; Bold=$false},
    @{Text="explaining the game play"; Bold=$false},
    @{Text=" instructions"; Bold=$false}
) $d

# --- "prevScore – stock label showing previous score" --------------------
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="prevScore"; Bold=$true},
    @{Text=" "; Bold=$true},
    @{Text=$dash; Bold=$true},
    @{Text=" "; Bold=$true},
    @{Text="stock label showing previous score"; Bold=$false}
) $d

# --- "btnPlay – stock button indicating “Play”" ---------------------------
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="btnPlay"; Bold=$true},
    @{Text=" $dash "; Bold=$true},
    @{Text="stock "; Bold=$false},
    @{Text="button indicating $ldq" + "Play$rdq"; Bold=$false}
) $d

# --- "btnQuit - stock button indicating "Quit”" ---------------------------
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="btnQuit"; Bold=$true},
    @{Text=" - "; Bold=$true},
    @{Text='stock button indicating "Quit' + $rdq; Bold=$false}
) $d

# --- "Other attributes:" ---------------------------------------------------
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="Other attributes:"; Bold=$false}
) $d

# --- "prevScore - integer indicating last score, passed into the class initializer and displayed on prevScore label"
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="prevScore"; Bold=$true},
    @{Text=" - integer indicating last score, passed into the class initializer and displayed on "; Bold=$false},
    @{Text="prevScore"; Bold=$false},
    @{Text=" label"; Bold=$false}
) $d

# --- "response - string representing the user's intentions. Set by the two buttons and read in the main function"
$idx = Add-RunsParagraphAtIndex $idx @(
    @{Text="response"; Bold=$true},
    @{Text=" - string representing the user's intentions. Set by the two buttons and read in the main function"; Bold=$false}
) $d

Write-Output "Inserted through paragraph index $idx"
